# Append two blank paragraphs, a new paragraph of text (with the
# spell-checker proofErr markers Word inserts around "linkingPages"),
# and a trailing blank paragraph to the end of the document body -
# mirroring a user typing this content after the last line then
# pressing Enter/Return a few times.

$d = $word.ActiveDocument

# Collapse a Range to the very end of the document's main story so the
# new content lands immediately before the closing sectPr, i.e. right
# after the last paragraph.
$endRange = $d.Content
$endRange.Collapse(0)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newContentXml = (
    "<w:p $wNs/>" +
    "<w:p $wNs/>" +
    "<w:p $wNs>" +
        "<w:r><w:t xml:space=`"preserve`">This Line of code was added while working with </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>linkingPages</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> branch.</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs/>"
)

$endRange.InsertXML($newContentXml)
